$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# This reverts an earlier edit to reference #72 (rows 17-18 of the sheet) that had
# rewritten the Abstract (column D) and Authors (column E) cells with slightly
# different (whitespace-padded / tag-stripped) text. Restore the original text so
# the previously-introduced shared strings become unused again.

$rightSingleQuote = [char]0x2019   # U+2019 RIGHT SINGLE QUOTATION MARK ( ' )
$alpha            = [char]0x03B1   # U+03B1 GREEK SMALL LETTER ALPHA ( a )
$thinSpace        = [char]0x2009   # U+2009 THIN SPACE

$abstract  = "Background`nid=`"Par1`">Since December 2019, 2019 novel coronavirus pneumonia emerged in Wuhan city and rapidly spread throughout China and even the world."
$abstract += "`n`n We sought to analyse the clinical characteristics and laboratory findings of some cases with 2019 novel coronavirus pneumonia ."
$abstract += "`n`n`nMethods`nid=`"Par2`">In this retrospective study, we extracted the data on 95 patients with laboratory-confirmed 2019 novel coronavirus pneumonia in Wuhan Xinzhou District People"
$abstract += $rightSingleQuote
$abstract += "s Hospital from January 16th to February 25th, 2020. Cases were confirmed by real-time RT-PCR and abnormal radiologic findings."
$abstract += "`n`n Outcomes were followed up until March 2th, 2020.`nResults`nid=`"Par3`">Higher temperature, blood leukocyte count, neutrophil count, neutrophil percentage, C-reactive protein level, D-dimer level, alanine aminotransferase activity, aspartate aminotransferase activity, "
$abstract += $alpha
$abstract += " - hydroxybutyrate dehydrogenase activity, lactate dehydrogenase activity and creatine kinase activity were related to severe 2019 novel coronavirus pneumonia and composite endpoint, and so were lower lymphocyte count, lymphocyte percentage and total protein level."
$abstract += "`n`n Age below 40 or above 60"
$abstract += $thinSpace
$abstract += "years old, male, higher Creatinine level, and lower platelet count also seemed related to severe 2019 novel coronavirus pneumonia and composite endpoint, however the P values were greater than 0.05, which mean under the same condition studies of larger samples are needed in the future."
$abstract += "`n`n`nConclusion`nid=`"Par4`">Multiple factors were related to severe 2019 novel coronavirus pneumonia and composite endpoint, and more related studies are needed in the future.`n`n`n"

$authors17 = "[Gemin%Zhang%NULL%0,    Jie%Zhang%945128911@qq.com%2,    Bowen%Wang%NULL%2,    Xionglin%Zhu%NULL%2,    Qiang%Wang%NULL%2,    Shiming%Qiu%NULL%2]"
$authors18 = "[Gemin%Zhang%NULL%0,    Jie%Zhang%945128911@qq.com%0,    Bowen%Wang%NULL%0,    Xionglin%Zhu%NULL%0,    Qiang%Wang%NULL%0,    Shiming%Qiu%NULL%0]"

$ws.Range("D17").Value = $abstract
$ws.Range("E17").Value = $authors17

$ws.Range("D18").Value = $abstract
$ws.Range("E18").Value = $authors18
